$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Case_4_120 loading_percent values (rows 2-25, columns B-E,G,I,J,L-O).
# Columns F, H, K stay 0; column A (index) is unchanged.

# Row 2
$ws.Cells.Item(2, 2).Value = 18.83343114763041
$ws.Cells.Item(2, 3).Value = 6.977096061341165
$ws.Cells.Item(2, 4).Value = 13.15180704825282
$ws.Cells.Item(2, 5).Value = 13.26763969818576
$ws.Cells.Item(2, 7).Value = 3.697519966356513
$ws.Cells.Item(2, 9).Value = 28.55227858295811
$ws.Cells.Item(2, 10).Value = 8.189554730806986
$ws.Cells.Item(2, 12).Value = 12.63843372616258
$ws.Cells.Item(2, 13).Value = 18.35963301309227
$ws.Cells.Item(2, 14).Value = 19.89982915868515
$ws.Cells.Item(2, 15).Value = 30.57403029773107

# Row 3
$ws.Cells.Item(3, 2).Value = 18.51092996293036
$ws.Cells.Item(3, 3).Value = 6.61960954460661
$ws.Cells.Item(3, 4).Value = 13.16736052531826
$ws.Cells.Item(3, 5).Value = 13.29850883145722
$ws.Cells.Item(3, 7).Value = 3.699902617033609
$ws.Cells.Item(3, 9).Value = 28.65156297293163
$ws.Cells.Item(3, 10).Value = 8.184686054880942
$ws.Cells.Item(3, 12).Value = 12.6466073225927
$ws.Cells.Item(3, 13).Value = 18.2976834503562
$ws.Cells.Item(3, 14).Value = 19.96199452755652
$ws.Cells.Item(3, 15).Value = 30.6401517173179

# Row 4
$ws.Cells.Item(4, 2).Value = 18.31329158645985
$ws.Cells.Item(4, 3).Value = 6.388569354170738
$ws.Cells.Item(4, 4).Value = 13.17905421553124
$ws.Cells.Item(4, 5).Value = 13.31861325287
$ws.Cells.Item(4, 7).Value = 3.701444170316042
$ws.Cells.Item(4, 9).Value = 28.71778242412224
$ws.Cells.Item(4, 10).Value = 8.18166812583511
$ws.Cells.Item(4, 12).Value = 12.65323051135943
$ws.Cells.Item(4, 13).Value = 18.26220719989518
$ws.Cells.Item(4, 14).Value = 20.00201850220244
$ws.Cells.Item(4, 15).Value = 30.68716754912372

# Row 5
$ws.Cells.Item(5, 2).Value = 18.23295099683034
$ws.Cells.Item(5, 3).Value = 6.291569463587981
$ws.Cells.Item(5, 4).Value = 13.18435869668625
$ws.Cells.Item(5, 5).Value = 13.32709594916721
$ws.Cells.Item(5, 7).Value = 3.702092193282865
$ws.Cells.Item(5, 9).Value = 28.74608808264416
$ws.Cells.Item(5, 10).Value = 8.180431156597766
$ws.Cells.Item(5, 12).Value = 12.65633370229695
$ws.Cells.Item(5, 13).Value = 18.24840428268873
$ws.Cells.Item(5, 14).Value = 20.01879621296351
$ws.Cells.Item(5, 15).Value = 30.70793735907537

# Row 6
$ws.Cells.Item(6, 2).Value = 18.21962554425965
$ws.Cells.Item(6, 3).Value = 6.275292392479128
$ws.Cells.Item(6, 4).Value = 13.18527207027461
$ws.Cells.Item(6, 5).Value = 13.32852202750157
$ws.Cells.Item(6, 7).Value = 3.702200996412865
$ws.Cells.Item(6, 9).Value = 28.75086793596123
$ws.Cells.Item(6, 10).Value = 8.180225328751961
$ws.Cells.Item(6, 12).Value = 12.65687341589309
$ws.Cells.Item(6, 13).Value = 18.2461520940109
$ws.Cells.Item(6, 14).Value = 20.02161042218464
$ws.Cells.Item(6, 15).Value = 30.71148334162839

# Row 7
$ws.Cells.Item(7, 2).Value = 18.31220714341778
$ws.Cells.Item(7, 3).Value = 6.38727263761006
$ws.Cells.Item(7, 4).Value = 13.17912357028911
$ws.Cells.Item(7, 5).Value = 13.31872647849835
$ws.Cells.Item(7, 7).Value = 3.701452829414943
$ws.Cells.Item(7, 9).Value = 28.71815881881603
$ws.Cells.Item(7, 10).Value = 8.181651472488552
$ws.Cells.Item(7, 12).Value = 12.65327072466063
$ws.Cells.Item(7, 13).Value = 18.26201838846197
$ws.Cells.Item(7, 14).Value = 20.00224287711166
$ws.Cells.Item(7, 15).Value = 30.68744114168086

# Row 8
$ws.Cells.Item(8, 2).Value = 18.72221582471442
$ws.Cells.Item(8, 3).Value = 6.856258970959401
$ws.Cells.Item(8, 4).Value = 13.15672505345237
$ws.Cells.Item(8, 5).Value = 13.27804494428261
$ws.Cells.Item(8, 7).Value = 3.698325230170232
$ws.Cells.Item(8, 9).Value = 28.58541950786331
$ws.Cells.Item(8, 10).Value = 8.187881957567496
$ws.Cells.Item(8, 12).Value = 12.64091939799592
$ws.Cells.Item(8, 13).Value = 18.3377474100549
$ws.Cells.Item(8, 14).Value = 19.92087976539213
$ws.Cells.Item(8, 15).Value = 30.59549531301657

# Row 9
$ws.Cells.Item(9, 2).Value = 19.52477817567895
$ws.Cells.Item(9, 3).Value = 7.682519224576851
$ws.Cells.Item(9, 4).Value = 13.12980263299086
$ws.Cells.Item(9, 5).Value = 13.20736938394323
$ws.Cells.Item(9, 7).Value = 3.692812732967951
$ws.Cells.Item(9, 9).Value = 28.36691174251039
$ws.Cells.Item(9, 10).Value = 8.199872202923718
$ws.Cells.Item(9, 12).Value = 12.62939713319671
$ws.Cells.Item(9, 13).Value = 18.50611025902711
$ws.Cells.Item(9, 14).Value = 19.77597618335007
$ws.Cells.Item(9, 15).Value = 30.46623559705564

# Row 10
$ws.Cells.Item(10, 2).Value = 20.10770696436656
$ws.Cells.Item(10, 3).Value = 8.230733827700053
$ws.Cells.Item(10, 4).Value = 13.12037004322059
$ws.Cells.Item(10, 5).Value = 13.16095171797371
$ws.Cells.Item(10, 7).Value = 3.689137038720172
$ws.Cells.Item(10, 9).Value = 28.2319481538627
$ws.Cells.Item(10, 10).Value = 8.208541603207566
$ws.Cells.Item(10, 12).Value = 12.62862622417773
$ws.Cells.Item(10, 13).Value = 18.64125870159048
$ws.Cells.Item(10, 14).Value = 19.67835720013046
$ws.Cells.Item(10, 15).Value = 30.40254883499279

# Row 11
$ws.Cells.Item(11, 2).Value = 20.37026737081575
$ws.Cells.Item(11, 3).Value = 8.467088343658686
$ws.Cells.Item(11, 4).Value = 13.11831971549927
$ws.Cells.Item(11, 5).Value = 13.14102252982238
$ws.Cells.Item(11, 7).Value = 3.687545289093014
$ws.Cells.Item(11, 9).Value = 28.17612323030873
$ws.Cells.Item(11, 10).Value = 8.212454974764297
$ws.Cells.Item(11, 12).Value = 12.62993471563773
$ws.Cells.Item(11, 13).Value = 18.70507753167107
$ws.Cells.Item(11, 14).Value = 19.63584882996004
$ws.Cells.Item(11, 15).Value = 30.38039608368004

# Row 12
$ws.Cells.Item(12, 2).Value = 20.46922342627746
$ws.Cells.Item(12, 3).Value = 8.554699926230201
$ws.Cells.Item(12, 4).Value = 13.11786474904628
$ws.Cells.Item(12, 5).Value = 13.13364583551739
$ws.Cells.Item(12, 7).Value = 3.686954022290878
$ws.Cells.Item(12, 9).Value = 28.15578656504951
$ws.Cells.Item(12, 10).Value = 8.213932460474025
$ws.Cells.Item(12, 12).Value = 12.6306676447977
$ws.Cells.Item(12, 13).Value = 18.72956690025342
$ws.Cells.Item(12, 14).Value = 19.62002369939533
$ws.Cells.Item(12, 15).Value = 30.3729895527885

# Row 13
$ws.Cells.Item(13, 2).Value = 20.44793380905476
$ws.Cells.Item(13, 3).Value = 8.535915521010178
$ws.Cells.Item(13, 4).Value = 13.11794845107326
$ws.Cells.Item(13, 5).Value = 13.13522698401924
$ws.Cells.Item(13, 7).Value = 3.687080851899096
$ws.Cells.Item(13, 9).Value = 28.16013068417341
$ws.Cells.Item(13, 10).Value = 8.213614455091065
$ws.Cells.Item(13, 12).Value = 12.63049925491093
$ws.Cells.Item(13, 13).Value = 18.72427855490248
$ws.Cells.Item(13, 14).Value = 19.62341985210523
$ws.Cells.Item(13, 15).Value = 30.37454097338571

# Row 14
$ws.Cells.Item(14, 2).Value = 20.37841848205073
$ws.Cells.Item(14, 3).Value = 8.474334176693082
$ws.Cells.Item(14, 4).Value = 13.11827584760024
$ws.Cells.Item(14, 5).Value = 13.14041223968901
$ws.Cells.Item(14, 7).Value = 3.687496415191259
$ws.Cells.Item(14, 9).Value = 28.1744340157108
$ws.Cells.Item(14, 10).Value = 8.212576617365704
$ws.Cells.Item(14, 12).Value = 12.62999026261518
$ws.Cells.Item(14, 13).Value = 18.70708590273969
$ws.Cells.Item(14, 14).Value = 19.63454144518546
$ws.Cells.Item(14, 15).Value = 30.37976705002763

# Row 15
$ws.Cells.Item(15, 2).Value = 20.33577444491506
$ws.Cells.Item(15, 3).Value = 8.436367164121915
$ws.Cells.Item(15, 4).Value = 13.11851822495736
$ws.Cells.Item(15, 5).Value = 13.14361048865205
$ws.Cells.Item(15, 7).Value = 3.687752454856486
$ws.Cells.Item(15, 9).Value = 28.18329985425577
$ws.Cells.Item(15, 10).Value = 8.211940334131063
$ws.Cells.Item(15, 12).Value = 12.62970937492333
$ws.Cells.Item(15, 13).Value = 18.69659649691399
$ws.Cells.Item(15, 14).Value = 19.64138911160795
$ws.Cells.Item(15, 15).Value = 30.38309613317479

# Row 16
$ws.Cells.Item(16, 2).Value = 20.09048726690601
$ws.Cells.Item(16, 3).Value = 8.215023745175031
$ws.Cells.Item(16, 4).Value = 13.12054906279079
$ws.Cells.Item(16, 5).Value = 13.16227796416099
$ws.Cells.Item(16, 7).Value = 3.689242674947526
$ws.Cells.Item(16, 9).Value = 28.23570870880697
$ws.Cells.Item(16, 10).Value = 8.208285239200647
$ws.Cells.Item(16, 12).Value = 12.62857400134535
$ws.Cells.Item(16, 13).Value = 18.63713389066463
$ws.Cells.Item(16, 14).Value = 19.68117333122081
$ws.Cells.Item(16, 15).Value = 30.40413393436908

# Row 17
$ws.Cells.Item(17, 2).Value = 19.93927332363981
$ws.Cells.Item(17, 3).Value = 8.075885770824385
$ws.Cells.Item(17, 4).Value = 13.12236829918606
$ws.Cells.Item(17, 5).Value = 13.17403333667958
$ws.Cells.Item(17, 7).Value = 3.690177412223603
$ws.Cells.Item(17, 9).Value = 28.26928803537751
$ws.Cells.Item(17, 10).Value = 8.206035137093586
$ws.Cells.Item(17, 12).Value = 12.6283016482237
$ws.Cells.Item(17, 13).Value = 18.60124493884834
$ws.Cells.Item(17, 14).Value = 19.7060651759576
$ws.Cells.Item(17, 15).Value = 30.41878771665051

# Row 18
$ws.Cells.Item(18, 2).Value = 19.85205808242076
$ws.Cells.Item(18, 3).Value = 7.994631456177827
$ws.Cells.Item(18, 4).Value = 13.123625625445
$ws.Cells.Item(18, 5).Value = 13.18090642202751
$ws.Cells.Item(18, 7).Value = 3.690722614052035
$ws.Cells.Item(18, 9).Value = 28.28912619949577
$ws.Cells.Item(18, 10).Value = 8.204738081503175
$ws.Cells.Item(18, 12).Value = 12.62830124278915
$ws.Cells.Item(18, 13).Value = 18.58082358820969
$ws.Cells.Item(18, 14).Value = 19.72056112013209
$ws.Cells.Item(18, 15).Value = 30.42785791230709

# Row 19
$ws.Cells.Item(19, 2).Value = 19.82249003518625
$ws.Cells.Item(19, 3).Value = 7.966910188362789
$ws.Cells.Item(19, 4).Value = 13.12408758589008
$ws.Cells.Item(19, 5).Value = 13.18325273317935
$ws.Cells.Item(19, 7).Value = 3.690908511236206
$ws.Cells.Item(19, 9).Value = 28.29593303673725
$ws.Cells.Item(19, 10).Value = 8.204298431811436
$ws.Cells.Item(19, 12).Value = 12.62832797450131
$ws.Cells.Item(19, 13).Value = 18.57394764418494
$ws.Cells.Item(19, 14).Value = 19.72549994618075
$ws.Cells.Item(19, 15).Value = 30.43103908336146

# Row 20
$ws.Cells.Item(20, 2).Value = 19.95539593853949
$ws.Cells.Item(20, 3).Value = 8.090824250375283
$ws.Cells.Item(20, 4).Value = 13.12215281086936
$ws.Cells.Item(20, 5).Value = 13.1727704000792
$ws.Cells.Item(20, 7).Value = 3.690077125254367
$ws.Cells.Item(20, 9).Value = 28.26565918923346
$ws.Cells.Item(20, 10).Value = 8.206274960318007
$ws.Cells.Item(20, 12).Value = 12.62831447801316
$ws.Cells.Item(20, 13).Value = 18.60504260229523
$ws.Cells.Item(20, 14).Value = 19.70339689884678
$ws.Cells.Item(20, 15).Value = 30.4171613663846

# Row 21
$ws.Cells.Item(21, 2).Value = 20.39885029438905
$ws.Cells.Item(21, 3).Value = 8.492473538669326
$ws.Cells.Item(21, 4).Value = 13.11817096592599
$ws.Cells.Item(21, 5).Value = 13.13888459347999
$ws.Cells.Item(21, 7).Value = 3.687374042695578
$ws.Cells.Item(21, 9).Value = 28.17021097121572
$ws.Cells.Item(21, 10).Value = 8.212881575946097
$ws.Cells.Item(21, 12).Value = 12.63013333181839
$ws.Cells.Item(21, 13).Value = 18.71212715934236
$ws.Cells.Item(21, 14).Value = 19.63126739485758
$ws.Cells.Item(21, 15).Value = 30.37820535530926

# Row 22
$ws.Cells.Item(22, 2).Value = 20.68589439044951
$ws.Cells.Item(22, 3).Value = 8.743948807370366
$ws.Cells.Item(22, 4).Value = 13.11744180311583
$ws.Cells.Item(22, 5).Value = 13.11772922963368
$ws.Cells.Item(22, 7).Value = 3.685674392192858
$ws.Cells.Item(22, 9).Value = 28.11251156106396
$ws.Cells.Item(22, 10).Value = 8.2171735734247
$ws.Cells.Item(22, 12).Value = 12.6327053927427
$ws.Cells.Item(22, 13).Value = 18.78398654710775
$ws.Cells.Item(22, 14).Value = 19.585710692892
$ws.Cells.Item(22, 15).Value = 30.35847108605423

# Row 23
$ws.Cells.Item(23, 2).Value = 20.53297757311235
$ws.Cells.Item(23, 3).Value = 8.61074512586832
$ws.Cells.Item(23, 4).Value = 13.11765985497339
$ws.Cells.Item(23, 5).Value = 13.12892975221495
$ws.Cells.Item(23, 7).Value = 3.686575419205951
$ws.Cells.Item(23, 9).Value = 28.14287782089177
$ws.Cells.Item(23, 10).Value = 8.214885233283816
$ws.Cells.Item(23, 12).Value = 12.63120646537094
$ws.Cells.Item(23, 13).Value = 18.74546709047811
$ws.Cells.Item(23, 14).Value = 19.60988061145649
$ws.Cells.Item(23, 15).Value = 30.36847926630844

# Row 24
$ws.Cells.Item(24, 2).Value = 19.94810777031021
$ws.Cells.Item(24, 3).Value = 8.084074493485703
$ws.Cells.Item(24, 4).Value = 13.12224957452952
$ws.Cells.Item(24, 5).Value = 13.17334101619413
$ws.Cells.Item(24, 7).Value = 3.690122440664175
$ws.Cells.Item(24, 9).Value = 28.2672981305357
$ws.Cells.Item(24, 10).Value = 8.206166546957959
$ws.Cells.Item(24, 12).Value = 12.62830819110999
$ws.Cells.Item(24, 13).Value = 18.60332501820892
$ws.Cells.Item(24, 14).Value = 19.7046026496428
$ws.Cells.Item(24, 15).Value = 30.41789462860909

# Row 25
$ws.Cells.Item(25, 2).Value = 19.30843287710019
$ws.Cells.Item(25, 3).Value = 7.469217447136494
$ws.Cells.Item(25, 4).Value = 13.13526704342888
$ws.Cells.Item(25, 5).Value = 13.2255188971444
$ws.Cells.Item(25, 7).Value = 3.694237979787997
$ws.Cells.Item(25, 9).Value = 28.42153966514568
$ws.Cells.Item(25, 10).Value = 8.196653001341481
$ws.Cells.Item(25, 12).Value = 12.62970937492333
$ws.Cells.Item(25, 13).Value = 18.45850293998463
$ws.Cells.Item(25, 14).Value = 19.8136173880764
$ws.Cells.Item(25, 15).Value = 30.37171933317479

Write-Output "Updated loading_percent values for rows 2-25"
